$wb = $excel.ActiveWorkbook

# Sheet "建物" (building) - property_category column (I) rows 2-6 were wrongly
# tagged "land"; fix to "building".
$wsBuilding = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 6; $r++) {
    $wsBuilding.Cells.Item($r, 9).Value = "building"
}

# Sheet "汽車" (car) - property_category column (H) row 2 was wrongly tagged
# "land"; fix to "car".
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Cells.Item(2, 8).Value = "car"
